$d = $word.ActiveDocument

function Merge-RunsBeforeSymbol([string]$searchText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $start = $r.Start
        $delRange = $d.Range($r.Start, $r.End)
        $delRange.Text = ""
        $insRange = $d.Range($start, $start)
        $insRange.InsertBefore($searchText)
    }
}

# Merge adjacent same-format runs in the letter->digit list items (no visible text change)
Merge-RunsBeforeSymbol("D, E, or F ")
Merge-RunsBeforeSymbol("G, H, or I ")
Merge-RunsBeforeSymbol("J, K, or L ")
Merge-RunsBeforeSymbol("M, N, or O ")

# Re-locate the split point inside "...numbers and print the result..." (between "p" and "rint")
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("numbers and p", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPos = $r.End

# Moves the document's _GoBack bookmark from the end of the document to the edit point,
# splitting the run there (mirrors Word's own "last edit" position tracking).
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Add a right-aligned header with the author's name
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdr.Range.InsertAfter("Josh Guerra")
$hdr.Range.Paragraphs(1).Style = "Header"
$hdr.Range.ParagraphFormat.Alignment = 2
